$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price/Volume cells in this sheet are textual (e.g. "23.422.61",
# "  +0.88%  "). Several of the new values look numeric to Excel and
# would otherwise get silently converted to a Number (losing exact
# formatting/precision, e.g. "6.640" -> 6.64 or "0.3787" -> 0.37869999999999998).
# Force every D/E cell we touch to keep a Text format so the literal string is preserved.

$textCells = @("D2","E2","E3","D4","E4","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E22","D23","E23","D24","E24","D25","E25","D26","E26","E27","D28","E28","D29","E29","D30","E30","D31","E31","D32","E32","D33","E33","D34","E34","D35","E35","D36","E36","D37","E37","D38","E38","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","E45","D46","E46","D47","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "23.453.84"
$ws.Range("E2").Value = "  +1.01%  "

$ws.Range("E3").Value = "  +2.20%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("D6").Value = "304.19"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").Value = "0.3787"
$ws.Range("E7").Value = "  +0.69%  "

$ws.Range("D8").Value = "52.26"
$ws.Range("E8").Value = "  -1.44%  "

$ws.Range("D9").Value = "0.3645"
$ws.Range("E9").Value = "  +1.09%  "

$ws.Range("D10").Value = "1.249"
$ws.Range("E10").Value = "  -1.04%  "

$ws.Range("D11").Value = "0.08107"
$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.10%  "

$ws.Range("D13").Value = "22.92"
$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("D14").Value = "6.640"
$ws.Range("E14").Value = "  +0.60%  "

$ws.Range("D15").Value = "0.00001253"
$ws.Range("E15").Value = "  +0.88%  "

$ws.Range("D16").Value = "7.284"

$ws.Range("D17").Value = "1.639.37"
$ws.Range("E17").Value = "  +2.32%  "

$ws.Range("D18").Value = "94.15"
$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("D19").Value = "0.06939"
$ws.Range("E19").Value = "  +0.30%  "

$ws.Range("D20").Value = "18.16"
$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("D21").Value = "6.546"

$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").Value = "23.457.24"
$ws.Range("E23").Value = "  +1.00%  "

$ws.Range("D24").Value = "12.86"
$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("D25").Value = "3.258"
$ws.Range("E25").Value = "  +6.03%  "

$ws.Range("D26").Value = "2.457"
$ws.Range("E26").Value = "  +1.95%  "

$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").Value = "149.35"
$ws.Range("E28").Value = "  -0.83%  "

$ws.Range("D29").Value = "5.303"
$ws.Range("E29").Value = "  +0.83%  "

$ws.Range("D30").Value = "135.75"
$ws.Range("E30").Value = "  +0.61%  "

$ws.Range("D31").Value = "2.314"
$ws.Range("E31").Value = "  -3.86%  "

$ws.Range("D32").Value = "1.822.17"
$ws.Range("E32").Value = "  +2.31%  "

$ws.Range("D33").Value = "6.897"
$ws.Range("E33").Value = "  +2.38%  "

$ws.Range("D34").Value = "11.00"
$ws.Range("E34").Value = "  +7.21%  "

$ws.Range("D35").Value = "0.9660"
$ws.Range("E35").Value = "  +1.88%  "

$ws.Range("D36").Value = "0.02875"
$ws.Range("E36").Value = "  +3.94%  "

$ws.Range("D37").Value = "6.274"
$ws.Range("E37").Value = "  +2.71%  "

$ws.Range("D38").Value = "0.2562"
$ws.Range("E38").Value = "  +1.99%  "

$ws.Range("D39").Value = "0.07291"
$ws.Range("E39").Value = "  -1.39%  "

$ws.Range("D40").Value = "0.08899"
$ws.Range("E40").Value = "  +1.67%  "

$ws.Range("D41").Value = "1.376"
$ws.Range("E41").Value = "  -1.51%  "

$ws.Range("D42").Value = "0.7128"
$ws.Range("E42").Value = "  +0.45%  "

$ws.Range("D43").Value = "16.47"
$ws.Range("E43").Value = "  +3.83%  "

$ws.Range("D44").Value = "12.58"
$ws.Range("E44").Value = "  +1.19%  "

$ws.Range("E45").Value = "  +0.59%  "

$ws.Range("D46").Value = "2.360"
$ws.Range("E46").Value = "  +1.53%  "

$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("D48").Value = "4.000"
$ws.Range("E48").Value = "  -0.16%  "

$ws.Range("D49").Value = "0.07999"
$ws.Range("E49").Value = "  +0.38%  "

$ws.Range("D50").Value = "1.224"
$ws.Range("E50").Value = "  +2.40%  "

$ws.Range("D51").Value = "127.62"
$ws.Range("E51").Value = "  -4.68%  "
